$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '23.699.05'
$ws.Range("E2").Value = '  +0.93%  '
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '1.656.76'
$ws.Range("E3").Value = '  +0.89%  '
$rng.Style = "Normal"

$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$ws.Range("E5").Value = '  +0.08%  '
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = '302.87'
$ws.Range("E6").Value = '  -0.12%  '
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = '0.3842'
$ws.Range("E7").Value = '  +0.53%  '
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '0.3608'
$ws.Range("E8").Value = '  +0.13%  '
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = '51.12'
$ws.Range("E9").Value = '  -1.64%  '
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '0.08196'
$ws.Range("E10").Value = '  -0.82%  '
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '1.230'
$ws.Range("E11").Value = '  -0.12%  '
$rng.Style = "Normal"

$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$ws.Range("E12").Value = '  +0.05%  '
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = '22.38'
$ws.Range("E13").Value = '  -0.65%  '
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '6.449'
$ws.Range("E14").Value = '  -0.14%  '
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = '7.442'
$ws.Range("E15").Value = '  +1.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = '0.00001222'
$ws.Range("E16").Value = '  -1.46%  '
$rng.Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = '1.651.20'
$ws.Range("E17").Value = '  +0.85%  '
$rng.Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = '97.54'
$ws.Range("E18").Value = '  +2.54%  '
$rng.Style = "Normal"

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = '0.07044'
$ws.Range("E19").Value = '  +1.04%  '
$rng.Style = "Normal"

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '6.794'
$ws.Range("E20").Value = '  +3.01%  '
$rng.Style = "Normal"

$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$ws.Range("E21").Value = '  +0.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$ws.Range("E23").Value = '  +1.69%  '
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '23.694.34'
$ws.Range("E24").Value = '  +0.93%  '
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '2.488'
$ws.Range("E25").Value = '  -1.90%  '
$rng.Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = '3.030'
$ws.Range("E26").Value = '  -1.32%  '
$rng.Style = "Normal"

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = '21.25'
$ws.Range("E27").Value = '  +0.35%  '
$rng.Style = "Normal"

$rng = $ws.Range("D28:E28")
$rng.NumberFormat = "@"
$ws.Range("D28").Value = '153.80'
$ws.Range("E28").Value = '  +1.25%  '
$rng.Style = "Normal"

$rng = $ws.Range("D29:E29")
$rng.NumberFormat = "@"
$ws.Range("D29").Value = '5.235'
$ws.Range("E29").Value = '  -0.78%  '
$rng.Style = "Normal"

$rng = $ws.Range("D30:E30")
$rng.NumberFormat = "@"
$ws.Range("D30").Value = '134.04'
$ws.Range("E30").Value = '  +0.46%  '
$rng.Style = "Normal"

$rng = $ws.Range("D31:E31")
$rng.NumberFormat = "@"
$ws.Range("D31").Value = '1.839.60'
$ws.Range("E31").Value = '  +1.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("D32:E32")
$rng.NumberFormat = "@"
$ws.Range("D32").Value = '7.136'
$ws.Range("E32").Value = '  +9.26%  '
$rng.Style = "Normal"

$rng = $ws.Range("D33:E33")
$rng.NumberFormat = "@"
$ws.Range("D33").Value = '2.251'
$ws.Range("E33").Value = '  +4.59%  '
$rng.Style = "Normal"

$rng = $ws.Range("D34:E34")
$rng.NumberFormat = "@"
$ws.Range("D34").Value = '12.01'
$ws.Range("E34").Value = '  +4.46%  '
$rng.Style = "Normal"

$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$ws.Range("E35").Value = '  -3.23%  '
$rng.Style = "Normal"

$rng = $ws.Range("D36:E36")
$rng.NumberFormat = "@"
$ws.Range("D36").Value = '0.02801'
$ws.Range("E36").Value = '  +0.94%  '
$rng.Style = "Normal"

$rng = $ws.Range("B37:E37")
$rng.NumberFormat = "@"
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").Value = '0.08814'
$ws.Range("E37").Value = '  +0.27%  '
$rng.Style = "Normal"

$rng = $ws.Range("B38:E38")
$rng.NumberFormat = "@"
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2501'
$ws.Range("E38").Value = '  -0.50%  '
$rng.Style = "Normal"

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = '6.083'
$ws.Range("E39").Value = '  +1.89%  '
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = '0.06984'
$ws.Range("E40").Value = '  -0.78%  '
$rng.Style = "Normal"

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = '13.02'
$ws.Range("E41").Value = '  +6.23%  '
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = '0.6981'
$ws.Range("E42").Value = '  -0.81%  '
$rng.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = '1.337'
$ws.Range("E43").Value = '  -0.76%  '
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = '15.97'
$ws.Range("E44").Value = '  +2.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = '0.6505'
$ws.Range("E45").Value = '  -0.20%  '
$rng.Style = "Normal"

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$ws.Range("E46").Value = '  +0.08%  '
$rng.Style = "Normal"

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = '2.295'
$ws.Range("E47").Value = '  +0.07%  '
$rng.Style = "Normal"

$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$ws.Range("E48").Value = '  -0.05%  '
$rng.Style = "Normal"

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = '0.07895'
$ws.Range("E49").Value = '  -1.10%  '
$rng.Style = "Normal"

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = '128.12'
$ws.Range("E50").Value = '  -0.65%  '
$rng.Style = "Normal"

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = '1.178'
$ws.Range("E51").Value = '  -1.24%  '
$rng.Style = "Normal"
